$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("material-age-velocity")

$data = @(
    @(18, "wwmd", "kiosk + hydrant", -0.7966854185018963, 0.006704448355048949, 0.008156122999893109, 5358.539281606674, -0.0145914976373844, -0.1960026948329951, -0.1480531118982936, -0.06934176270382052, -0.01000927001696754),
    @(19, "wwmd", "kiosk + hydrant", -0.6625668846702162, 0.00154224852862009, 0.0009395606420866699, 7634.565146446228, -0.01029578055671602, -0.02213087515193186, -0.3651462109107909, -0.03409693106586862, -0.1276630981074167),
    @(16, "wwmd", "kiosk + hydrant", -0.6572141270400695, 0.0105418661906856, 0.007352559039659772, 16508.92272257805, -0.01000529444249022, -0.01844996715188763, -0.04490408910009593, -0.09613977806487756, -0.01000286903281968)
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
